# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.202.18"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "1.788.55"

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'226.13"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "'32.34"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").Value = "'0.0690"
$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("E11").Value = "  +0.78%  "

$ws.Range("D12").Value = "2.046.16"
$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("D13").Value = "'11.13"
$ws.Range("E13").Value = "  -2.27%  "

$ws.Range("D14").Value = "1.791.34"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").Value = "'0.625"
$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").Value = "34.176.98"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").Value = "0.0₃0807"
$ws.Range("E19").Value = "  +3.09%  "

$ws.Range("D20").Value = "'246.01"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("E21").Value = "  +0.63%  "

$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("E23").Value = "  +1.73%  "

$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").Value = "'161.93"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "'7.18"
$ws.Range("E26").Value = "  -0.45%  "

$ws.Range("D27").Value = "'16.33"
$ws.Range("E27").Value = "  +0.21%  "

$ws.Range("E28").Value = "  +0.79%  "

$ws.Range("E29").Value = "  +0.22%  "

$ws.Range("D30").Value = "'0.0521"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("E31").Value = "  -0.83%  "

$ws.Range("D32").Value = "'3.76"
$ws.Range("E32").Value = "  +2.71%  "

$ws.Range("D33").Value = "'3.75"
$ws.Range("E33").Value = "  +3.86%  "

$ws.Range("E34").Value = "  -2.05%  "

$ws.Range("D35").Value = "1.441.68"
$ws.Range("E35").Value = "  +1.97%  "

$ws.Range("E36").Value = "  +8.38%  "

$ws.Range("D37").Value = "'0.666"
$ws.Range("E37").Value = "  +2.89%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0191"
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.05"
$ws.Range("E39").Value = "  +1.00%  "

$ws.Range("D40").Value = "'82.14"
$ws.Range("E40").Value = "  +1.69%  "

$ws.Range("E41").Value = "  +1.32%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'13.86"
$ws.Range("E43").Value = "  +3.04%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.72"
$ws.Range("E44").Value = "  +0.83%  "

$ws.Range("D45").Value = "'0.0519"
$ws.Range("E45").Value = "  +2.31%  "

$ws.Range("D46").Value = "'6.10"
$ws.Range("E46").Value = "  +0.84%  "

$ws.Range("E47").Value = "  +0.64%  "

$ws.Range("D48").Value = "1.946.39"
$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("D49").Value = "'105.46"

$ws.Range("E50").Value = "  +0.17%  "

$ws.Range("E51").Value = "  -6.59%  "
